$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column E: header + values mirroring the "Id" column (1..6)
$ws.Range("E1").Value = "Nova Propriedade"
$ws.Range("E2").Value = 1
$ws.Range("E3").Value = 2
$ws.Range("E4").Value = 3
$ws.Range("E5").Value = 4
$ws.Range("E6").Value = 5
$ws.Range("E7").Value = 6

# Size column E to fit its contents (best-fit width)
$ws.Columns("E").ColumnWidth = 16.5

# Move / leave the selection on the last filled cell of the new column
$ws.Range("E7").Select()
